# Removed Saving role and committeeOf outside
# Only persistent data saved outside is password now

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "Role" (E) and "Committee Of" (F) columns - header and data -
# leaving only Name, Email, Faculty, Password as persisted columns.
$ws.Range("E1:F12").ClearContents()

# Reflect the resulting selection left after selecting column F (whole
# column) as the last user interaction before saving.
$ws.Columns.Item(6).Select()
